# Update Excel data - 2024-11-22 05:36:27
# Applies refreshed crypto market data to all three sheets of the workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws3 = $wb.Worksheets.Item("Summary")

$ws1.Cells.Item(2,3).Value = [double]"99049"
$ws1.Cells.Item(2,4).Value = [double]"1959195618520"
$ws1.Cells.Item(2,5).Value = [double]"112637342948"
$ws1.Cells.Item(2,6).Value = [double]"1.23955"
$ws1.Cells.Item(3,3).Value = [double]"3387.32"
$ws1.Cells.Item(3,4).Value = [double]"408077072603"
$ws1.Cells.Item(3,5).Value = [double]"56006401370"
$ws1.Cells.Item(3,6).Value = [double]"8.157629999999999"
$ws1.Cells.Item(4,4).Value = [double]"130946260507"
$ws1.Cells.Item(4,5).Value = [double]"118643630134"
$ws1.Cells.Item(4,6).Value = [double]"-0.04242"
$ws1.Cells.Item(5,3).Value = [double]"262.45"
$ws1.Cells.Item(5,4).Value = [double]"124537701624"
$ws1.Cells.Item(5,5).Value = [double]"14996288932"
$ws1.Cells.Item(5,6).Value = [double]"8.74217"
$ws1.Cells.Item(6,3).Value = [double]"635.4299999999999"
$ws1.Cells.Item(6,4).Value = [double]"92774471217"
$ws1.Cells.Item(6,5).Value = [double]"2471223191"
$ws1.Cells.Item(6,6).Value = [double]"3.69055"
$ws1.Cells.Item(7,4).Value = [double]"79043620118"
$ws1.Cells.Item(7,5).Value = [double]"17986070520"
$ws1.Cells.Item(7,6).Value = [double]"24.44778"
$ws1.Cells.Item(8,3).Value = [double]"0.396036"
$ws1.Cells.Item(8,4).Value = [double]"58263803542"
$ws1.Cells.Item(8,5).Value = [double]"9818155855"
$ws1.Cells.Item(8,6).Value = [double]"2.05484"
$ws1.Cells.Item(9,4).Value = [double]"38333818068"
$ws1.Cells.Item(9,5).Value = [double]"10390110683"
$ws1.Cells.Item(9,6).Value = [double]"-0.05674"
$ws1.Cells.Item(10,3).Value = [double]"3386.1"
$ws1.Cells.Item(10,4).Value = [double]"33167306357"
$ws1.Cells.Item(10,5).Value = [double]"148210338"
$ws1.Cells.Item(10,6).Value = [double]"8.32367"
$ws1.Cells.Item(11,3).Value = [double]"0.884086"
$ws1.Cells.Item(11,4).Value = [double]"31689028378"
$ws1.Cells.Item(11,5).Value = [double]"3620403635"
$ws1.Cells.Item(11,6).Value = [double]"12.00104"
$ws1.Cells.Item(12,3).Value = [double]"0.20032"
$ws1.Cells.Item(12,4).Value = [double]"17312612565"
$ws1.Cells.Item(12,5).Value = [double]"1068387702"
$ws1.Cells.Item(12,6).Value = [double]"1.32711"
$ws1.Cells.Item(13,3).Value = [double]"36.33"
$ws1.Cells.Item(13,4).Value = [double]"14861570205"
$ws1.Cells.Item(13,5).Value = [double]"1050610416"
$ws1.Cells.Item(13,6).Value = [double]"6.52267"
$ws1.Cells.Item(14,4).Value = [double]"14719052047"
$ws1.Cells.Item(14,5).Value = [double]"1610749228"
$ws1.Cells.Item(14,6).Value = [double]"2.89691"
$ws1.Cells.Item(15,3).Value = [double]"4011.57"
$ws1.Cells.Item(15,4).Value = [double]"14487656085"
$ws1.Cells.Item(15,5).Value = [double]"157676278"
$ws1.Cells.Item(15,6).Value = [double]"8.152979999999999"
$ws1.Cells.Item(16,3).Value = [double]"99016"
$ws1.Cells.Item(16,4).Value = [double]"14437890938"
$ws1.Cells.Item(16,5).Value = [double]"854187820"
$ws1.Cells.Item(16,6).Value = [double]"1.62661"
$ws1.Cells.Item(17,4).Value = [double]"14174714359"
$ws1.Cells.Item(17,5).Value = [double]"638662056"
$ws1.Cells.Item(17,6).Value = [double]"3.18693"
$ws1.Cells.Item(18,3).Value = [double]"3.62"
$ws1.Cells.Item(18,4).Value = [double]"10296343269"
$ws1.Cells.Item(18,5).Value = [double]"2081734581"
$ws1.Cells.Item(18,6).Value = [double]"0.5196"
$ws1.Cells.Item(19,3).Value = [double]"497.01"
$ws1.Cells.Item(19,4).Value = [double]"9820919497"
$ws1.Cells.Item(19,5).Value = [double]"1936718320"
$ws1.Cells.Item(19,6).Value = [double]"-3.23961"
$ws1.Cells.Item(20,3).Value = [double]"3386.67"
$ws1.Cells.Item(20,4).Value = [double]"9647766380"
$ws1.Cells.Item(20,5).Value = [double]"1272082578"
$ws1.Cells.Item(20,6).Value = [double]"8.148099999999999"
$ws1.Cells.Item(21,3).Value = [double]"15.29"
$ws1.Cells.Item(21,4).Value = [double]"9589028913"
$ws1.Cells.Item(21,5).Value = [double]"1250243694"
$ws1.Cells.Item(21,6).Value = [double]"4.19627"
$ws1.Cells.Item(22,3).Value = [double]"6.22"
$ws1.Cells.Item(22,4).Value = [double]"8950792803"
$ws1.Cells.Item(22,5).Value = [double]"828809348"
$ws1.Cells.Item(22,6).Value = [double]"8.66395"
$ws1.Cells.Item(23,3).Value = [double]"2.129e-05"
$ws1.Cells.Item(23,4).Value = [double]"8947854713"
$ws1.Cells.Item(23,5).Value = [double]"6816364706"
$ws1.Cells.Item(23,6).Value = [double]"9.08339"
$ws1.Cells.Item(24,3).Value = [double]"0.283467"
$ws1.Cells.Item(24,4).Value = [double]"8476827049"
$ws1.Cells.Item(24,5).Value = [double]"2299577623"
$ws1.Cells.Item(24,6).Value = [double]"18.22379"
$ws1.Cells.Item(25,3).Value = [double]"8.800000000000001"
$ws1.Cells.Item(25,4).Value = [double]"8141138033"
$ws1.Cells.Item(25,5).Value = [double]"3450901"
$ws1.Cells.Item(25,6).Value = [double]"3.34373"
$ws1.Cells.Item(26,3).Value = [double]"5.8"
$ws1.Cells.Item(26,4).Value = [double]"7070338866"
$ws1.Cells.Item(26,5).Value = [double]"1011815621"
$ws1.Cells.Item(26,6).Value = [double]"4.29611"
$ws1.Cells.Item(27,3).Value = [double]"90.54000000000001"
$ws1.Cells.Item(27,4).Value = [double]"6809385663"
$ws1.Cells.Item(27,5).Value = [double]"1413117460"
$ws1.Cells.Item(27,6).Value = [double]"4.08485"
$ws1.Cells.Item(28,3).Value = [double]"12.13"
$ws1.Cells.Item(28,4).Value = [double]"6470864453"
$ws1.Cells.Item(28,5).Value = [double]"865197468"
$ws1.Cells.Item(28,6).Value = [double]"3.80178"
$ws1.Cells.Item(29,3).Value = [double]"3568.39"
$ws1.Cells.Item(29,4).Value = [double]"6207201257"
$ws1.Cells.Item(29,5).Value = [double]"104283760"
$ws1.Cells.Item(29,6).Value = [double]"8.14161"
$ws1.Cells.Item(30,3).Value = [double]"9.41"
$ws1.Cells.Item(30,4).Value = [double]"5652647882"
$ws1.Cells.Item(30,5).Value = [double]"861384256"
$ws1.Cells.Item(30,6).Value = [double]"5.93018"
$ws1.Cells.Item(31,3).Value = [double]"0.200382"
$ws1.Cells.Item(31,4).Value = [double]"5446042457"
$ws1.Cells.Item(31,5).Value = [double]"126658464"
$ws1.Cells.Item(31,6).Value = [double]"14.08539"
$ws1.Cells.Item(32,3).Value = [double]"1.002"
$ws1.Cells.Item(32,4).Value = [double]"5226148199"
$ws1.Cells.Item(32,5).Value = [double]"16537152"
$ws1.Cells.Item(32,6).Value = [double]"0.00528"
$ws1.Cells.Item(33,3).Value = [double]"0.133661"
$ws1.Cells.Item(33,4).Value = [double]"5102145955"
$ws1.Cells.Item(33,5).Value = [double]"895575078"
$ws1.Cells.Item(33,6).Value = [double]"5.98247"
$ws1.Cells.Item(34,3).Value = [double]"9.630000000000001"
$ws1.Cells.Item(34,4).Value = [double]"4568754468"
$ws1.Cells.Item(34,5).Value = [double]"274366712"
$ws1.Cells.Item(34,6).Value = [double]"5.80122"
$ws1.Cells.Item(35,3).Value = [double]"27.99"
$ws1.Cells.Item(35,4).Value = [double]"4188436495"
$ws1.Cells.Item(35,5).Value = [double]"882892358"
$ws1.Cells.Item(35,6).Value = [double]"5.20897"
$ws1.Cells.Item(36,3).Value = [double]"5.195e-05"
$ws1.Cells.Item(36,4).Value = [double]"3901989396"
$ws1.Cells.Item(36,5).Value = [double]"1682777811"
$ws1.Cells.Item(36,6).Value = [double]"2.20426"
$ws1.Cells.Item(39,3).Value = [double]"0.472771"
$ws1.Cells.Item(39,4).Value = [double]"3768901866"
$ws1.Cells.Item(39,5).Value = [double]"491287376"
$ws1.Cells.Item(39,6).Value = [double]"7.38356"
$ws1.Cells.Item(40,3).Value = [double]"507.56"
$ws1.Cells.Item(40,4).Value = [double]"3757496611"
$ws1.Cells.Item(40,5).Value = [double]"284629438"
$ws1.Cells.Item(40,6).Value = [double]"3.05338"
$ws1.Cells.Item(41,3).Value = [double]"1.002"
$ws1.Cells.Item(41,4).Value = [double]"3687424253"
$ws1.Cells.Item(41,5).Value = [double]"224417787"
$ws1.Cells.Item(41,6).Value = [double]"-0.0735"
$ws1.Cells.Item(42,3).Value = [double]"24.82"
$ws1.Cells.Item(42,4).Value = [double]"3576586417"
$ws1.Cells.Item(42,5).Value = [double]"33415891"
$ws1.Cells.Item(42,6).Value = [double]"2.58434"
$ws1.Cells.Item(43,3).Value = [double]"1.001"
$ws1.Cells.Item(43,4).Value = [double]"3443652111"
$ws1.Cells.Item(43,5).Value = [double]"154533119"
$ws1.Cells.Item(43,6).Value = [double]"-0.07159"
$ws1.Cells.Item(44,3).Value = [double]"3.38"
$ws1.Cells.Item(44,4).Value = [double]"3386517689"
$ws1.Cells.Item(44,5).Value = [double]"1281177864"
$ws1.Cells.Item(44,6).Value = [double]"6.17607"
$ws1.Cells.Item(45,3).Value = [double]"3.71"
$ws1.Cells.Item(45,4).Value = [double]"3344314078"
$ws1.Cells.Item(45,5).Value = [double]"301495698"
$ws1.Cells.Item(45,6).Value = [double]"2.17889"
$ws1.Cells.Item(46,4).Value = [double]"3339864365"
$ws1.Cells.Item(46,5).Value = [double]"485294902"
$ws1.Cells.Item(46,6).Value = [double]"2.29024"
$ws1.Cells.Item(47,3).Value = [double]"0.78912"
$ws1.Cells.Item(47,4).Value = [double]"3234863626"
$ws1.Cells.Item(47,5).Value = [double]"1667684113"
$ws1.Cells.Item(47,6).Value = [double]"12.79568"
$ws1.Cells.Item(48,3).Value = [double]"161.42"
$ws1.Cells.Item(48,4).Value = [double]"2970687075"
$ws1.Cells.Item(48,5).Value = [double]"86854554"
$ws1.Cells.Item(48,6).Value = [double]"-1.16649"
$ws1.Cells.Item(49,3).Value = [double]"1.96"
$ws1.Cells.Item(49,4).Value = [double]"2941098154"
$ws1.Cells.Item(49,5).Value = [double]"367168437"
$ws1.Cells.Item(49,6).Value = [double]"0.98808"
$ws1.Cells.Item(50,3).Value = [double]"0.847303"
$ws1.Cells.Item(50,4).Value = [double]"2838053749"
$ws1.Cells.Item(50,5).Value = [double]"185133243"
$ws1.Cells.Item(50,6).Value = [double]"15.30567"
$ws1.Cells.Item(51,3).Value = [double]"4.7"
$ws1.Cells.Item(51,4).Value = [double]"2819630872"
$ws1.Cells.Item(51,5).Value = [double]"581494444"
$ws1.Cells.Item(51,6).Value = [double]"6.70952"
# Rows 37 and 38 swapped order (Kaspa now ranks above Render) with updated figures
$ws1.Cells.Item(37,1).Value = "Kaspa"
$ws1.Cells.Item(37,2).Value = "kas"
$ws1.Cells.Item(37,3).Value = [double]"0.151353"
$ws1.Cells.Item(37,4).Value = [double]"3817573608"
$ws1.Cells.Item(37,5).Value = [double]"151839983"
$ws1.Cells.Item(37,6).Value = [double]"-0.79223"

$ws1.Cells.Item(38,1).Value = "Render"
$ws1.Cells.Item(38,2).Value = "render"
$ws1.Cells.Item(38,3).Value = [double]"7.38"
$ws1.Cells.Item(38,4).Value = [double]"3817395639"
$ws1.Cells.Item(38,5).Value = [double]"434850443"
$ws1.Cells.Item(38,6).Value = [double]"-0.12395"

# "Top 5 by Market Cap" sheet - refreshed market capitalization figures
$ws2.Cells.Item(2,2).Value = [double]"1959195618520"
$ws2.Cells.Item(3,2).Value = [double]"408077072603"
$ws2.Cells.Item(4,2).Value = [double]"130946260507"
$ws2.Cells.Item(5,2).Value = [double]"124537701624"
$ws2.Cells.Item(6,2).Value = [double]"92774471217"

# "Summary" sheet - refreshed headline metrics (kept as text, matching source formatting)
$ws3.Range("B2").NumberFormat = "@"
$ws3.Cells.Item(2,2).Value = '$4363.16'
$ws3.Cells.Item(3,2).Value = "XRP (24.45%)"
$ws3.Cells.Item(4,2).Value = "Bitcoin Cash (-3.24%)"
